# Scheduled runner: refresh market-price-derived figures (Universalis data pull)
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 15152767
$ws.Range("I19").Value = 652.3333
$ws.Range("J19").Value = 25642692
$ws.Range("K19").Value = 652.3333
$ws.Range("L19").Value = 25642692
$ws.Range("M19").Value = -477.3333
$ws.Range("N19").Value = -25643042
$ws.Range("H40").Value = 6466.579
$ws.Range("I40").Value = 3860
$ws.Range("K40").Value = 3860
$ws.Range("M40").Value = -3685
$ws.Range("H92").Value = 397.9
$ws.Range("I92").Value = 397.9
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 397.9
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 850.1
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 3496.1667
$ws.Range("I98").Value = 3142.6
$ws.Range("K98").Value = 3142.6
$ws.Range("M98").Value = -1644.6
$ws.Range("H99").Value = 84375.25
$ws.Range("I99").Value = 334
$ws.Range("K99").Value = 1002
$ws.Range("M99").Value = 496
$ws.Range("H122").Value = 3496.1667
$ws.Range("I122").Value = 3142.6
$ws.Range("K122").Value = 9427.799999999999
$ws.Range("M122").Value = -6977.799999999999
$ws.Range("H127").Value = 3232.652
$ws.Range("I127").Value = 941.0909
$ws.Range("K127").Value = 2823.2727
$ws.Range("M127").Value = 2136.7273
$ws.Range("H137").Value = 10055.787
$ws.Range("I137").Value = 1066.7916
$ws.Range("K137").Value = 3200.3748
$ws.Range("M137").Value = -650.3748000000001
$ws.Range("H138").Value = 5626
$ws.Range("I138").Value = 3196.5
$ws.Range("J138").Value = 6165.8887
$ws.Range("K138").Value = 9589.5
$ws.Range("L138").Value = 18497.6661
$ws.Range("M138").Value = -4449.5
$ws.Range("N138").Value = -28777.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 836.587
$ws.Range("I32").Value = 839.3333
$ws.Range("J32").Value = 788.8
$ws.Range("K32").Value = 839.3333
$ws.Range("L32").Value = 788.8
$ws.Range("M32").Value = -552.3333
$ws.Range("N32").Value = -1362.8
$ws.Range("H45").Value = 60921.94
$ws.Range("I45").Value = 68511.53
$ws.Range("K45").Value = 68511.53
$ws.Range("M45").Value = -68134.53
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 20000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 20000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -19594
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 20000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -18596
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2314.9167
$ws.Range("I20").Value = 1797.5
$ws.Range("J20").Value = 2832.3333
$ws.Range("K20").Value = 1797.5
$ws.Range("L20").Value = 2832.3333
$ws.Range("M20").Value = -1550.5
$ws.Range("N20").Value = -3326.3333
$ws.Range("H86").Value = 588161.4399999999
$ws.Range("I86").Value = 851658.1
$ws.Range("J86").Value = 2613.2222
$ws.Range("K86").Value = 851658.1
$ws.Range("L86").Value = 2613.2222
$ws.Range("M86").Value = -850535.1
$ws.Range("N86").Value = -4859.2222
$ws.Range("H89").Value = 588161.4399999999
$ws.Range("I89").Value = 851658.1
$ws.Range("J89").Value = 2613.2222
$ws.Range("K89").Value = 4258290.5
$ws.Range("L89").Value = 13066.111
$ws.Range("M89").Value = -4252674.5
$ws.Range("N89").Value = -24298.111
$ws.Range("H94").Value = 1287.7142
$ws.Range("I94").Value = 995.6667
$ws.Range("K94").Value = 995.6667
$ws.Range("M94").Value = -544.6667
$ws.Range("H134").Value = 44611.89
$ws.Range("J134").Value = 89587.766
$ws.Range("L134").Value = 268763.298
$ws.Range("N134").Value = -273833.298

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1030.1578
$ws.Range("I16").Value = 970.9286
$ws.Range("K16").Value = 970.9286
$ws.Range("M16").Value = -683.9286
$ws.Range("H31").Value = 581200.0600000001
$ws.Range("I31").Value = 997988.4
$ws.Range("K31").Value = 997988.4
$ws.Range("M31").Value = -997693.4
$ws.Range("H34").Value = 581200.0600000001
$ws.Range("I34").Value = 997988.4
$ws.Range("K34").Value = 997988.4
$ws.Range("M34").Value = -997786.4
$ws.Range("H68").Value = 79929.664
$ws.Range("I68").Value = 79944.5
$ws.Range("K68").Value = 79944.5
$ws.Range("M68").Value = -79195.5
$ws.Range("H71").Value = 79929.664
$ws.Range("I71").Value = 79944.5
$ws.Range("K71").Value = 239833.5
$ws.Range("M71").Value = -236089.5
$ws.Range("H113").Value = 1030.1578
$ws.Range("I113").Value = 970.9286
$ws.Range("K113").Value = 970.9286
$ws.Range("M113").Value = 1199.0714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4000.5806
$ws.Range("I14").Value = 4000.5806
$ws.Range("K14").Value = 12001.7418
$ws.Range("M14").Value = -11828.7418
$ws.Range("H68").Value = 2102.6667
$ws.Range("J68").Value = 2086.6667
$ws.Range("L68").Value = 6260.000100000001
$ws.Range("N68").Value = -7882.000100000001
$ws.Range("H71").Value = 2102.6667
$ws.Range("J71").Value = 2086.6667
$ws.Range("L71").Value = 18780.0003
$ws.Range("N71").Value = -26892.0003
$ws.Range("H141").Value = 4318.75
$ws.Range("I141").Value = 3166.6667
$ws.Range("K141").Value = 9500.000100000001
$ws.Range("M141").Value = -4320.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1543011.4
$ws.Range("I80").Value = 1001464.2
$ws.Range("K80").Value = 1001464.2
$ws.Range("M80").Value = -1000466.2
$ws.Range("H83").Value = 1543011.4
$ws.Range("I83").Value = 1001464.2
$ws.Range("K83").Value = 5007321
$ws.Range("M83").Value = -5002329
$ws.Range("H97").Value = 539.5
$ws.Range("I97").Value = 602.28
$ws.Range("K97").Value = 602.28
$ws.Range("M97").Value = -106.28
$ws.Range("H113").Value = 48122530
$ws.Range("I113").Value = 1494663.1
$ws.Range("J113").Value = 71436456
$ws.Range("K113").Value = 1494663.1
$ws.Range("L113").Value = 71436456
$ws.Range("M113").Value = -1492493.1
$ws.Range("N113").Value = -71440796
$ws.Range("H132").Value = 52172.297
$ws.Range("I132").Value = 21959.883
$ws.Range("K132").Value = 65879.649
$ws.Range("M132").Value = -63349.649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 561632.2
$ws.Range("I40").Value = 631123.7
$ws.Range("K40").Value = 631123.7
$ws.Range("M40").Value = -630987.7
$ws.Range("H82").Value = 1041.5
$ws.Range("I82").Value = 487.5
$ws.Range("J82").Value = 2149.5
$ws.Range("K82").Value = 487.5
$ws.Range("L82").Value = 2149.5
$ws.Range("M82").Value = -126.5
$ws.Range("N82").Value = -2871.5
$ws.Range("H85").Value = 1041.5
$ws.Range("I85").Value = 487.5
$ws.Range("J85").Value = 2149.5
$ws.Range("K85").Value = 487.5
$ws.Range("L85").Value = 2149.5
$ws.Range("M85").Value = 760.5
$ws.Range("N85").Value = -4645.5
$ws.Range("H122").Value = 554952.6
$ws.Range("I122").Value = 3806.7
$ws.Range("K122").Value = 11420.1
$ws.Range("M122").Value = -8970.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H81").Value = 3408.5557
$ws.Range("I81").Value = 2953.8572
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 5907.7144
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -4846.7144
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 3408.5557
$ws.Range("I84").Value = 2953.8572
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 29538.572
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -24234.572
$ws.Range("N84").Value = -60608
$ws.Range("H113").Value = 1866.25
$ws.Range("I113").Value = 2248.75
$ws.Range("J113").Value = 1101.25
$ws.Range("K113").Value = 6746.25
$ws.Range("L113").Value = 3303.75
$ws.Range("M113").Value = -4576.25
$ws.Range("N113").Value = -7643.75
$ws.Range("H126").Value = 8357.75
$ws.Range("I126").Value = 1885.909
$ws.Range("K126").Value = 5657.727000000001
$ws.Range("M126").Value = -3187.727000000001
$ws.Range("H136").Value = 569304
$ws.Range("I136").Value = 608237.9399999999
$ws.Range("K136").Value = 1824713.82
$ws.Range("M136").Value = -1822163.82
